# Applies the diff: updates rows 899-918 (date + reshuffled product/price
# values for "Disco" on 2025-06-01) and appends 21 new rows (1819-1839) of
# "Hiper Libertad" purchases for the same date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "YYYY-MM-DD HH:MM:SS"
$newDate = 45809

# --- Update existing rows 899-918 -----------------------------------------
# row => (B value, C value)
$updates = @{
    899 = @("Arroz", 2350)
    900 = @("Galletas de agua", 1522)
    901 = @("Galletas dulces", 2850)
    902 = @("Batata", 1399)
    903 = @("Papa", 969)
    904 = @("Harina", 950)
    905 = @("Dulce de Leche", 3500)
    906 = @("Azucar", 1200)
    907 = @("Pollo", 8497.5)
    908 = @("Leche", 1950)
    909 = @("Salame", 6376)
    910 = @("Gaseosa", 4200)
    911 = @("Manteca", 2400)
    912 = @("Aceite", 4500)
    913 = @("Yerba", 3300)
    914 = @("Sal fina", 1200)
    915 = @("Café", 12450)
    916 = @("Vinagre", 1500)
    917 = @("Cerveza", 3550)
    918 = @("Mayonesa", 1350)
}

foreach ($row in 899..918) {
    $ws.Cells.Item($row, 1).Value = $newDate
    $vals = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
}

# --- Append new rows 1819-1839 ---------------------------------------------
# Each entry: (B value, C value)
$newRows = @(
    @("Batata", 1699),
    @("Galletas dulces", 2149),
    @("Azucar", 1099),
    @("Galletas de agua", 1529),
    @("Papa", 1099),
    @("Arroz", 1999),
    @("Pan", 989),
    @("Queso crema", 4795),
    @("Dulce de Leche", 2342),
    @("Manzana", 2599),
    @("Huevo", 2205),
    @("Leche", 1766),
    @("Salame", 3590),
    @("Yogur", 3375),
    @("Pollo", 8290),
    @("Sal fina", 1550),
    @("Yerba", 3490),
    @("Café", 10489),
    @("Vinagre", 1115),
    @("Mayonesa", 1120),
    @("Manteca", 1580)
)

$startRow = 1819
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $entry = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $newDate
    $ws.Cells.Item($r, 1).NumberFormat = $dateFormat
    $ws.Cells.Item($r, 2).Value = $entry[0]
    $ws.Cells.Item($r, 3).Value = $entry[1]
    $ws.Cells.Item($r, 4).Value = "Hiper Libertad"
}
